$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H62").Value = 5373.4
$ws.Range("I62").Value = 5373.4
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5373.4
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4749.4

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H65").Value = 5373.4
$ws.Range("I65").Value = 5373.4
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 26867
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -23747

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H70").Value = 1903.5
$ws.Range("I70").Value = 1885.5
$ws.Range("J70").Value = 1948.5
$ws.Range("K70").Value = 5656.5
$ws.Range("L70").Value = 5845.5
$ws.Range("M70").Value = -5386.5
$ws.Range("N70").Value = -6385.5

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H73").Value = 1903.5
$ws.Range("I73").Value = 1885.5
$ws.Range("J73").Value = 1948.5
$ws.Range("K73").Value = 5656.5
$ws.Range("L73").Value = 5845.5
$ws.Range("M73").Value = -4720.5
$ws.Range("N73").Value = -7717.5

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H106").Value = 2534.075
$ws.Range("I106").Value = 1656.303
$ws.Range("J106").Value = 6672.143
$ws.Range("K106").Value = 1656.303
$ws.Range("L106").Value = 6672.143
$ws.Range("M106").Value = -1025.303
$ws.Range("N106").Value = -7934.143

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H137").Value = 2155.16
$ws.Range("I137").Value = 1896.0834
$ws.Range("J137").Value = 2543.775
$ws.Range("K137").Value = 5688.2502
$ws.Range("L137").Value = 7631.325000000001
$ws.Range("M137").Value = -3138.2502
$ws.Range("N137").Value = -12731.325

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 1679.4546
$ws.Range("I2").Value = 1682.2433
$ws.Range("J2").Value = 1664.7142
$ws.Range("K2").Value = 1682.2433
$ws.Range("L2").Value = 1664.7142
$ws.Range("M2").Value = -1569.2433
$ws.Range("N2").Value = -1890.7142

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 2631.111
$ws.Range("I32").Value = 1767.537
$ws.Range("J32").Value = 7812.5557
$ws.Range("K32").Value = 1767.537
$ws.Range("L32").Value = 7812.5557
$ws.Range("M32").Value = -1480.537
$ws.Range("N32").Value = -8386.555700000001

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 10000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -9727

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H74").Value = 1539.5938
$ws.Range("I74").Value = 1437.6428
$ws.Range("J74").Value = 2253.25
$ws.Range("K74").Value = 1437.6428
$ws.Range("L74").Value = 2253.25
$ws.Range("M74").Value = -563.6428000000001
$ws.Range("N74").Value = -4001.25

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H77").Value = 1539.5938
$ws.Range("I77").Value = 1437.6428
$ws.Range("J77").Value = 2253.25
$ws.Range("K77").Value = 7188.214
$ws.Range("L77").Value = 11266.25
$ws.Range("M77").Value = -2820.214
$ws.Range("N77").Value = -20002.25

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H97").Value = 3331.5386
$ws.Range("I97").Value = 3310.348
$ws.Range("J97").Value = 3494
$ws.Range("K97").Value = 3310.348
$ws.Range("L97").Value = 3494
$ws.Range("M97").Value = -2814.348
$ws.Range("N97").Value = -4486

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H116").Value = 1679.4546
$ws.Range("I116").Value = 1682.2433
$ws.Range("J116").Value = 1664.7142
$ws.Range("K116").Value = 1682.2433
$ws.Range("L116").Value = 1664.7142
$ws.Range("M116").Value = 611.7566999999999
$ws.Range("N116").Value = -6252.7142

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 1679.4546
$ws.Range("I3").Value = 1682.2433
$ws.Range("J3").Value = 1664.7142
$ws.Range("K3").Value = 1682.2433
$ws.Range("L3").Value = 1664.7142
$ws.Range("M3").Value = -1568.2433
$ws.Range("N3").Value = -1892.7142

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H105").Value = 2289.0857
$ws.Range("I105").Value = 2106.64
$ws.Range("J105").Value = 2745.2
$ws.Range("K105").Value = 2106.64
$ws.Range("L105").Value = 2745.2
$ws.Range("M105").Value = -359.6399999999999
$ws.Range("N105").Value = -6239.2

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 6012.2607
$ws.Range("I134").Value = 3905.0588
$ws.Range("J134").Value = 11982.667
$ws.Range("K134").Value = 11715.1764
$ws.Range("L134").Value = 35948.001
$ws.Range("M134").Value = -9180.1764
$ws.Range("N134").Value = -41018.001

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 2301.4807
$ws.Range("I31").Value = 2165.75
$ws.Range("J31").Value = 2518.65
$ws.Range("K31").Value = 2165.75
$ws.Range("L31").Value = 2518.65
$ws.Range("M31").Value = -1870.75
$ws.Range("N31").Value = -3108.65

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 2301.4807
$ws.Range("I34").Value = 2165.75
$ws.Range("J34").Value = 2518.65
$ws.Range("K34").Value = 2165.75
$ws.Range("L34").Value = 2518.65
$ws.Range("M34").Value = -1963.75
$ws.Range("N34").Value = -2922.65

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H93").Value = 5904.6665
$ws.Range("I93").Value = 5904.6665
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 5904.6665
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -4032.6665

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H134").Value = 4856.868
$ws.Range("I134").Value = 4380.5386
$ws.Range("J134").Value = 6183.7856
$ws.Range("K134").Value = 13141.6158
$ws.Range("L134").Value = 18551.3568
$ws.Range("M134").Value = -10606.6158
$ws.Range("N134").Value = -23621.3568

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H45").Value = 2752.8333
$ws.Range("I45").Value = 3740
$ws.Range("J45").Value = 778.5
$ws.Range("K45").Value = 11220
$ws.Range("L45").Value = 2335.5
$ws.Range("M45").Value = -10688
$ws.Range("N45").Value = -3399.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H68").Value = 1834.6
$ws.Range("I68").Value = 1182
$ws.Range("J68").Value = 2813.5
$ws.Range("K68").Value = 3546
$ws.Range("L68").Value = 8440.5
$ws.Range("M68").Value = -2735
$ws.Range("N68").Value = -10062.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H71").Value = 1834.6
$ws.Range("I71").Value = 1182
$ws.Range("J71").Value = 2813.5
$ws.Range("K71").Value = 10638
$ws.Range("L71").Value = 25321.5
$ws.Range("M71").Value = -6582
$ws.Range("N71").Value = -33433.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H98").Value = 962.8182
$ws.Range("I98").Value = 1069.5
$ws.Range("J98").Value = 939.1111
$ws.Range("K98").Value = 3208.5
$ws.Range("L98").Value = 2817.3333
$ws.Range("M98").Value = -1710.5
$ws.Range("N98").Value = -5813.3333

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H132").Value = 1280.375
$ws.Range("I132").Value = 1280.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11523.375
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8993.375

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H64").Value = 59971
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 59971
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 59971
$ws.Range("N64").Value = -60467

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H67").Value = 59971
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 59971
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 59971
$ws.Range("N67").Value = -61687

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H99").Value = 16871.285
$ws.Range("I99").Value = 16871.285
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 16871.285
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -14625.285

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H122").Value = 2755.6365
$ws.Range("I122").Value = 2176.1428
$ws.Range("J122").Value = 3769.75
$ws.Range("K122").Value = 6528.428400000001
$ws.Range("L122").Value = 11309.25
$ws.Range("M122").Value = -4078.428400000001
$ws.Range("N122").Value = -16209.25

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H122").Value = 6798.7
$ws.Range("I122").Value = 5462.6665
$ws.Range("J122").Value = 7371.2856
$ws.Range("K122").Value = 16387.9995
$ws.Range("L122").Value = 22113.8568
$ws.Range("M122").Value = -13937.9995
$ws.Range("N122").Value = -27013.8568

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H45").Value = 12511.25
$ws.Range("I45").Value = 19000
$ws.Range("J45").Value = 10348.333
$ws.Range("K45").Value = 19000
$ws.Range("L45").Value = 10348.333
$ws.Range("M45").Value = -18509
$ws.Range("N45").Value = -11330.333

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H99").Value = 70000
$ws.Range("I99").Value = 70000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 70000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -67005

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H100").Value = 834.8421
$ws.Range("I100").Value = 567.8125
$ws.Range("J100").Value = 2259
$ws.Range("K100").Value = 1135.625
$ws.Range("L100").Value = 4518
$ws.Range("M100").Value = -594.625
$ws.Range("N100").Value = -5600

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 3258.7568
$ws.Range("I132").Value = 3258.7568
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9776.270400000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7246.270400000001

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H136").Value = 26190.576
$ws.Range("I136").Value = 26174.285
$ws.Range("J136").Value = 26281.8
$ws.Range("K136").Value = 78522.855
$ws.Range("L136").Value = 78845.39999999999
$ws.Range("M136").Value = -75972.855
$ws.Range("N136").Value = -83945.39999999999
